$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" date (column C) from 45184 to 45185 for rows 2 through 20
for ($row = 2; $row -le 20; $row++) {
    $ws.Cells.Item($row, 3).Value = 45185
}
